$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions): update "想去人数" (F column) for two events
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 305
$wsExhibit.Range("F6").Value = 64

# Sheet "全部类型" (all types): update "想去人数" (F column) for the same two events
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 305
$wsAll.Range("F7").Value = 64
